$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# (single-decimal-point numeric strings) - force them to stay Text so the
# stored value matches the literal string from the refreshed crypto feed.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated cell values from the crypto-list refresh
$ws.Range('D2').Value = '27.532.95'
$ws.Range('E2').Value = '  +5.16%  '
$ws.Range('D3').Value = '1.724.49'
$ws.Range('E3').Value = '  +4.06%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '226.39'
$ws.Range('E5').Value = '  +3.47%  '
$ws.Range('D6').Value = '0.5379'
$ws.Range('E6').Value = '  +2.63%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.2691'
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('D9').Value = '0.06623'
$ws.Range('E9').Value = '  +4.16%  '
$ws.Range('E10').Value = '  +5.46%  '
$ws.Range('D11').Value = '0.07757'
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = '4.651'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = '1.735.55'
$ws.Range('E13').Value = '  +3.62%  '
$ws.Range('D14').Value = '1.962.60'
$ws.Range('E14').Value = '  +4.08%  '
$ws.Range('D15').Value = '0.5898'
$ws.Range('E15').Value = '  +4.86%  '
$ws.Range('D16').Value = '0.0₅8301'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D18').Value = '27.549.59'
$ws.Range('E18').Value = '  +5.26%  '
$ws.Range('D19').Value = '225.26'
$ws.Range('E19').Value = '  +17.38%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '4.744'
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('D23').Value = '6.122'
$ws.Range('E23').Value = '  +2.55%  '
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '148.16'
$ws.Range('E25').Value = '  +1.64%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = '0.1233'
$ws.Range('E26').Value = '  +2.72%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '1.685'
$ws.Range('E27').Value = '  +10.55%  '
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('D29').Value = '16.76'
$ws.Range('E29').Value = '  +4.80%  '
$ws.Range('D30').Value = '0.05586'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('E31').Value = '  +2.78%  '
$ws.Range('D32').Value = '3.584'
$ws.Range('E32').Value = '  +2.95%  '
$ws.Range('D33').Value = '3.476'
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').Value = '1.670'
$ws.Range('E34').Value = '  +6.44%  '
$ws.Range('D35').Value = '0.9638'
$ws.Range('E35').Value = '  +1.12%  '
$ws.Range('D36').Value = '2.448'
$ws.Range('E36').Value = '  +1.87%  '
$ws.Range('D37').Value = '2.816'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').Value = '0.5960'
$ws.Range('E38').Value = '  +4.24%  '
$ws.Range('E39').Value = '  +3.50%  '
$ws.Range('D40').Value = '5.899'
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').Value = '0.8607'
$ws.Range('E41').Value = '  +3.14%  '
$ws.Range('D42').Value = '1.062.37'
$ws.Range('E42').Value = '  +2.80%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '101.77'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('D45').Value = '1.867.84'
$ws.Range('E45').Value = '  +3.98%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '59.13'
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '8.243'
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.4431'
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').Value = '1.005'
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.05285'
$ws.Range('E50').Value = '  +0.91%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '1.470'
$ws.Range('E51').Value = '  +5.35%  '

Write-Host "Applied 104 cell updates"
